$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il1rn"
$ws.Cells.Item(2, 3).Value = "Il1r1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 305.026516
$ws.Cells.Item(2, 8).Value = 915.079548
$ws.Cells.Item(2, 9).Value = 0.9998851412135495
$ws.Cells.Item(2, 10).Value = 0.9998851412135495
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.175804666666667
$ws.Cells.Item(2, 14).Value = 21.527414
$ws.Cells.Item(2, 15).Value = 0.1618789277039843
$ws.Cells.Item(2, 16).Value = 0.1618789277039842
$ws.Cells.Item(2, 17).Value = 2188.810696969875
$ws.Cells.Item(2, 18).Value = 19699.29627272887
$ws.Cells.Item(2, 19).Value = 0.1618603344867963
$ws.Cells.Item(2, 20).Value = 0.1618603344867963

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il1rn"
$ws.Cells.Item(3, 3).Value = "Il1r1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 305.026516
$ws.Cells.Item(3, 8).Value = 915.079548
$ws.Cells.Item(3, 9).Value = 0.9998851412135495
$ws.Cells.Item(3, 10).Value = 0.9998851412135495
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 24.92162533333333
$ws.Cells.Item(3, 14).Value = 74.764876
$ws.Cells.Item(3, 15).Value = 0.5622067730383848
$ws.Cells.Item(3, 16).Value = 0.5622067730383847
$ws.Cells.Item(3, 17).Value = 7601.756548484006
$ws.Cells.Item(3, 18).Value = 68415.80893635606
$ws.Cells.Item(3, 19).Value = 0.5621421986506995
$ws.Cells.Item(3, 20).Value = 0.5621421986506994

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il1rn"
$ws.Cells.Item(4, 3).Value = "Il1r1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 305.026516
$ws.Cells.Item(4, 8).Value = 915.079548
$ws.Cells.Item(4, 9).Value = 0.9998851412135495
$ws.Cells.Item(4, 10).Value = 0.9998851412135495
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 12.23078966666667
$ws.Cells.Item(4, 14).Value = 36.692369
$ws.Cells.Item(4, 15).Value = 0.2759142992576309
$ws.Cells.Item(4, 16).Value = 0.2759142992576308
$ws.Cells.Item(4, 17).Value = 3730.715159952135
$ws.Cells.Item(4, 18).Value = 33576.43643956921
$ws.Cells.Item(4, 19).Value = 0.2758826080760538
$ws.Cells.Item(4, 20).Value = 0.2758826080760537

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Il1rn"
$ws.Cells.Item(5, 3).Value = "Il1r1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.035039
$ws.Cells.Item(5, 8).Value = 0.105117
$ws.Cells.Item(5, 9).Value = 0.0001148587864504919
$ws.Cells.Item(5, 10).Value = 0.0001148587864504919
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 7.175804666666667
$ws.Cells.Item(5, 14).Value = 21.527414
$ws.Cells.Item(5, 15).Value = 0.1618789277039843
$ws.Cells.Item(5, 16).Value = 0.1618789277039842
$ws.Cells.Item(5, 17).Value = 0.2514330197153333
$ws.Cells.Item(5, 18).Value = 2.262897177438
$ws.Cells.Item(5, 19).Value = 0.00001859321718798655
$ws.Cells.Item(5, 20).Value = 0.00001859321718798655

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Il1rn"
$ws.Cells.Item(6, 3).Value = "Il1r1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.035039
$ws.Cells.Item(6, 8).Value = 0.105117
$ws.Cells.Item(6, 9).Value = 0.0001148587864504919
$ws.Cells.Item(6, 10).Value = 0.0001148587864504919
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 24.92162533333333
$ws.Cells.Item(6, 14).Value = 74.764876
$ws.Cells.Item(6, 15).Value = 0.5622067730383848
$ws.Cells.Item(6, 16).Value = 0.5622067730383847
$ws.Cells.Item(6, 17).Value = 0.8732288300546668
$ws.Cells.Item(6, 18).Value = 7.859059470492
$ws.Cells.Item(6, 19).Value = 0.00006457438768543604
$ws.Cells.Item(6, 20).Value = 0.00006457438768543602

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Il1rn"
$ws.Cells.Item(7, 3).Value = "Il1r1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.035039
$ws.Cells.Item(7, 8).Value = 0.105117
$ws.Cells.Item(7, 9).Value = 0.0001148587864504919
$ws.Cells.Item(7, 10).Value = 0.0001148587864504919
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 12.23078966666667
$ws.Cells.Item(7, 14).Value = 36.692369
$ws.Cells.Item(7, 15).Value = 0.2759142992576309
$ws.Cells.Item(7, 16).Value = 0.2759142992576308
$ws.Cells.Item(7, 17).Value = 0.4285546391303333
$ws.Cells.Item(7, 18).Value = 3.856991752173
$ws.Cells.Item(7, 19).Value = 0.00003169118157706935
$ws.Cells.Item(7, 20).Value = 0.00003169118157706934
